$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("production_rates")
$ws.Range("A1").Value = "test"
